# "Update 5 Desember 2021"
# Insert a new "Lingkungan" column before the "Id"/"Status"/... block (old
# column E), reusing the header style of the "Waktu Hadir" column, resize a
# few columns, move the selection, and set the print paper size.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new column ------------------------------------------------
# A plain EntireColumn.Insert() at E copies formatting from the LEFT
# neighbour (column D), which is not the style we want on the new header
# cell. Instead, insert a throw-away column after the last column (after H,
# to the right of "Waktu Hadir" which carries the style we want to reuse),
# which copies format from its left neighbour (the "Waktu Hadir" header).
# Then cut that freshly-formatted column and insert it (shifting cells
# right) at column E - this relocates the correctly-styled cell without
# disturbing rows that have no data out at column E onward.
$ws.Columns("H:H").Insert()
$ws.Columns("H:H").Cut()
$ws.Columns("E:E").Insert()

# Remove any stray formatted-but-empty cells the insert may have produced
# below the header row in the new column.
$ws.Range("E2:E6").Clear()

# --- New header ------------------------------------------------------------
$ws.Range("E1").Value = "Lingkungan"

# --- Column widths -----------------------------------------------------
$ws.Columns("A:A").ColumnWidth = 5.333333333333333
$ws.Columns("B:B").ColumnWidth = 11.666666666666666
$ws.Columns("C:C").ColumnWidth = 31.333333333333332
$ws.Columns("D:D").ColumnWidth = 13.333333333333334
$ws.Columns("E:E").ColumnWidth = 10.5
$ws.Columns("F:F").ColumnWidth = 13
$ws.Columns("G:G").ColumnWidth = 27.666666666666668
$ws.Columns("H:H").ColumnWidth = 11.5

# --- View / print settings -------------------------------------------------
[void]$ws.Range("C6").Select()
$ws.PageSetup.PaperSize = 9

$wb.Save()
